# Fix wrong height value for the default-styled chrome export row (DPI=100)
# and update the active selection to reflect the last-used cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the Height value (column F, row 15) from 365 to 353.
# Dependent formulas in H14, H16, H17, H18 (which divide by $F$15)
# will recalculate automatically.
$ws.Range("F15").Value = 353

# Reflect the new active cell selection on the sheet.
$ws.Activate()
$ws.Range("I15").Select()
